$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row398 = New-Object "object[,]" 1,28
$row398[0,0] = 3724261
$row398[0,1] = 'Poland I Liga'
$row398[0,2] = 'Poland I Liga'
$row398[0,3] = 44661.41666666666
$row398[0,4] = 'GKS Jastrzebie'
$row398[0,5] = 'Chrobry Glogow'
$row398[0,6] = 4
$row398[0,7] = 2
$row398[0,8] = 'H'
$row398[0,9] = 2.625
$row398[0,10] = 3
$row398[0,11] = 2.625
$row398[0,12] = 2.8
$row398[0,13] = 2.875
$row398[0,14] = 2.55
$row398[0,15] = 0
$row398[0,16] = 2.025
$row398[0,17] = 1.825
$row398[0,18] = 2
$row398[0,19] = 1.95
$row398[0,20] = 1.9
$row398[0,21] = 1.8
$row398[0,22] = -1
$row398[0,23] = -1
$row398[0,24] = 1.025
$row398[0,25] = -1
$row398[0,26] = 0.95
$row398[0,27] = -1
$ws.Range("B398:AC398").Value = $row398

$row399 = New-Object "object[,]" 1,28
$row399[0,0] = 3724262
$row399[0,1] = 'Poland I Liga'
$row399[0,2] = 'Poland I Liga'
$row399[0,3] = 44661.41666666666
$row399[0,4] = 'Widzew Lodz'
$row399[0,5] = 'Skra Czestochowa'
$row399[0,6] = 2
$row399[0,7] = 2
$row399[0,8] = 'D'
$row399[0,9] = 1.75
$row399[0,10] = 3.5
$row399[0,11] = 4
$row399[0,12] = 1.615
$row399[0,13] = 3.5
$row399[0,14] = 4.75
$row399[0,15] = -0.75
$row399[0,16] = 1.825
$row399[0,17] = 1.975
$row399[0,18] = 2.25
$row399[0,19] = 1.9
$row399[0,20] = 1.9
$row399[0,21] = -1
$row399[0,22] = 2.5
$row399[0,23] = -1
$row399[0,24] = -1
$row399[0,25] = 0.9750000000000001
$row399[0,26] = 0.8999999999999999
$row399[0,27] = -1
$ws.Range("B399:AC399").Value = $row399

$row510 = New-Object "object[,]" 1,28
$row510[0,0] = 5138951
$row510[0,1] = 'Poland I Liga'
$row510[0,2] = 'Poland I Liga'
$row510[0,3] = 44784.54166666666
$row510[0,4] = 'Termalica BB Nieciecza'
$row510[0,5] = 'Chrobry Glogow'
$row510[0,6] = 3
$row510[0,7] = 1
$row510[0,8] = 'H'
$row510[0,9] = 1.666
$row510[0,10] = 3.6
$row510[0,11] = 4.333
$row510[0,12] = 1.666
$row510[0,13] = 3.5
$row510[0,14] = 4.5
$row510[0,15] = -0.75
$row510[0,16] = 1.875
$row510[0,17] = 1.925
$row510[0,18] = 2.5
$row510[0,19] = 1.875
$row510[0,20] = 1.925
$row510[0,21] = 0.6659999999999999
$row510[0,22] = -1
$row510[0,23] = -1
$row510[0,24] = 0.875
$row510[0,25] = -1
$row510[0,26] = 0.875
$row510[0,27] = -1
$ws.Range("B510:AC510").Value = $row510

$row511 = New-Object "object[,]" 1,28
$row511[0,0] = 5140780
$row511[0,1] = 'Poland I Liga'
$row511[0,2] = 'Poland I Liga'
$row511[0,3] = 44784.54166666666
$row511[0,4] = 'Chojniczanka Chojnice'
$row511[0,5] = 'Stal Rzeszow'
$row511[0,6] = 1
$row511[0,7] = 1
$row511[0,8] = 'D'
$row511[0,9] = 2.55
$row511[0,10] = 3.2
$row511[0,11] = 2.55
$row511[0,12] = 2.625
$row511[0,13] = 3.4
$row511[0,14] = 2.375
$row511[0,15] = 0
$row511[0,16] = 1.975
$row511[0,17] = 1.825
$row511[0,18] = 2.75
$row511[0,19] = 1.95
$row511[0,20] = 1.9
$row511[0,21] = -1
$row511[0,22] = 2.4
$row511[0,23] = -1
$row511[0,24] = 0
$row511[0,25] = -0
$row511[0,26] = -1
$row511[0,27] = 0.8999999999999999
$ws.Range("B511:AC511").Value = $row511

$row587 = New-Object "object[,]" 1,28
$row587[0,0] = 5139001
$row587[0,1] = 'Poland I Liga'
$row587[0,2] = 'Poland I Liga'
$row587[0,3] = 44849.41666666666
$row587[0,4] = 'MKS Puszcza Niepolomice'
$row587[0,5] = 'GKS Katowice'
$row587[0,6] = 1
$row587[0,7] = 1
$row587[0,8] = 'D'
$row587[0,9] = 2.15
$row587[0,10] = 3.3
$row587[0,11] = 3
$row587[0,12] = 2.45
$row587[0,13] = 3.1
$row587[0,14] = 2.75
$row587[0,15] = 0
$row587[0,16] = 1.825
$row587[0,17] = 1.975
$row587[0,18] = 2.25
$row587[0,19] = 1.95
$row587[0,20] = 1.85
$row587[0,21] = -1
$row587[0,22] = 2.1
$row587[0,23] = -1
$row587[0,24] = 0
$row587[0,25] = -0
$row587[0,26] = -0.5
$row587[0,27] = 0.425
$ws.Range("B587:AC587").Value = $row587

$row588 = New-Object "object[,]" 1,28
$row588[0,0] = 5138999
$row588[0,1] = 'Poland I Liga'
$row588[0,2] = 'Poland I Liga'
$row588[0,3] = 44849.41666666666
$row588[0,4] = 'Gornik Leczna'
$row588[0,5] = 'Skra Czestochowa'
$row588[0,6] = 3
$row588[0,7] = 0
$row588[0,8] = 'H'
$row588[0,9] = 2.05
$row588[0,10] = 3.3
$row588[0,11] = 3.25
$row588[0,12] = 1.95
$row588[0,13] = 3.3
$row588[0,14] = 3.5
$row588[0,15] = -0.5
$row588[0,16] = 2
$row588[0,17] = 1.8
$row588[0,18] = 2.5
$row588[0,19] = 2
$row588[0,20] = 1.8
$row588[0,21] = 0.95
$row588[0,22] = -1
$row588[0,23] = -1
$row588[0,24] = 1
$row588[0,25] = -1
$row588[0,26] = 1
$row588[0,27] = -1
$ws.Range("B588:AC588").Value = $row588

$row628 = New-Object "object[,]" 1,28
$row628[0,0] = 5139019
$row628[0,1] = 'Poland I Liga'
$row628[0,2] = 'Poland I Liga'
$row628[0,3] = 44878.58333333334
$row628[0,4] = 'Gornik Leczna'
$row628[0,5] = 'MKS Puszcza Niepolomice'
$row628[0,6] = 2
$row628[0,7] = 2
$row628[0,8] = 'D'
$row628[0,9] = 2.25
$row628[0,10] = 3.3
$row628[0,11] = 3
$row628[0,12] = 2.5
$row628[0,13] = 3.2
$row628[0,14] = 2.625
$row628[0,15] = 0
$row628[0,16] = 1.875
$row628[0,17] = 1.925
$row628[0,18] = 2.5
$row628[0,19] = 2
$row628[0,20] = 1.8
$row628[0,21] = -1
$row628[0,22] = 2.2
$row628[0,23] = -1
$row628[0,24] = 0
$row628[0,25] = -0
$row628[0,26] = 1
$row628[0,27] = -1
$ws.Range("B628:AC628").Value = $row628

$row629 = New-Object "object[,]" 1,28
$row629[0,0] = 5139023
$row629[0,1] = 'Poland I Liga'
$row629[0,2] = 'Poland I Liga'
$row629[0,3] = 44878.58333333334
$row629[0,4] = 'Resovia Rzeszow'
$row629[0,5] = 'Zaglebie Sosnowiec'
$row629[0,6] = 2
$row629[0,7] = 2
$row629[0,8] = 'D'
$row629[0,9] = 2.375
$row629[0,10] = 3.25
$row629[0,11] = 2.625
$row629[0,12] = 2.8
$row629[0,13] = 3.3
$row629[0,14] = 2.25
$row629[0,15] = 0.25
$row629[0,16] = 1.8
$row629[0,17] = 2
$row629[0,18] = 2.5
$row629[0,19] = 1.95
$row629[0,20] = 1.85
$row629[0,21] = -1
$row629[0,22] = 2.3
$row629[0,23] = -1
$row629[0,24] = 0.4
$row629[0,25] = -0.5
$row629[0,26] = 0.95
$row629[0,27] = -1
$ws.Range("B629:AC629").Value = $row629

$row679 = New-Object "object[,]" 1,28
$row679[0,0] = 5139053
$row679[0,1] = 'Poland I Liga'
$row679[0,2] = 'Poland I Liga'
$row679[0,3] = 45004.58333333334
$row679[0,4] = 'Chrobry Glogow'
$row679[0,5] = 'Zaglebie Sosnowiec'
$row679[0,6] = 0
$row679[0,7] = 0
$row679[0,8] = 'D'
$row679[0,9] = 2.45
$row679[0,10] = 3.2
$row679[0,11] = 2.55
$row679[0,12] = 2.7
$row679[0,13] = 3.2
$row679[0,14] = 2.375
$row679[0,15] = 0
$row679[0,16] = 2.05
$row679[0,17] = 1.75
$row679[0,18] = 2.25
$row679[0,19] = 1.875
$row679[0,20] = 1.925
$row679[0,21] = -1
$row679[0,22] = 2.2
$row679[0,23] = -1
$row679[0,24] = 0
$row679[0,25] = -0
$row679[0,26] = -1
$row679[0,27] = 0.925
$ws.Range("B679:AC679").Value = $row679

$row680 = New-Object "object[,]" 1,28
$row680[0,0] = 5140743
$row680[0,1] = 'Poland I Liga'
$row680[0,2] = 'Poland I Liga'
$row680[0,3] = 45004.58333333334
$row680[0,4] = 'Stal Rzeszow'
$row680[0,5] = 'Termalica BB Nieciecza'
$row680[0,6] = 2
$row680[0,7] = 2
$row680[0,8] = 'D'
$row680[0,9] = 3
$row680[0,10] = 3.3
$row680[0,11] = 2.2
$row680[0,12] = 2.9
$row680[0,13] = 3.3
$row680[0,14] = 2.25
$row680[0,15] = 0.25
$row680[0,16] = 1.825
$row680[0,17] = 1.975
$row680[0,18] = 2.5
$row680[0,19] = 1.95
$row680[0,20] = 1.85
$row680[0,21] = -1
$row680[0,22] = 2.3
$row680[0,23] = -1
$row680[0,24] = 0.4125
$row680[0,25] = -0.5
$row680[0,26] = 0.95
$row680[0,27] = -1
$ws.Range("B680:AC680").Value = $row680

$row682 = New-Object "object[,]" 1,28
$row682[0,0] = 5139054
$row682[0,1] = 'Poland I Liga'
$row682[0,2] = 'Poland I Liga'
$row682[0,3] = 45004.58333333334
$row682[0,4] = 'GKS Tychy 71'
$row682[0,5] = 'Sandecja Nowy Sacz'
$row682[0,6] = 2
$row682[0,7] = 3
$row682[0,8] = 'A'
$row682[0,9] = 2.15
$row682[0,10] = 3.2
$row682[0,11] = 3.1
$row682[0,12] = 2.375
$row682[0,13] = 3
$row682[0,14] = 3
$row682[0,15] = -0.25
$row682[0,16] = 2.025
$row682[0,17] = 1.775
$row682[0,18] = 2.25
$row682[0,19] = 1.975
$row682[0,20] = 1.825
$row682[0,21] = -1
$row682[0,22] = -1
$row682[0,23] = 2
$row682[0,24] = -1
$row682[0,25] = 0.7749999999999999
$row682[0,26] = 0.9750000000000001
$row682[0,27] = -1
$ws.Range("B682:AC682").Value = $row682

$row848 = New-Object "object[,]" 1,28
$row848[0,0] = 6803738
$row848[0,1] = 'Poland I Liga'
$row848[0,2] = 'Poland I Liga'
$row848[0,3] = 45191.54166666666
$row848[0,4] = 'Podbeskidzie Bielsko Biala'
$row848[0,5] = 'Gornik Leczna'
$row848[0,6] = 1
$row848[0,7] = 1
$row848[0,8] = 'D'
$row848[0,9] = 1.85
$row848[0,10] = 3.5
$row848[0,11] = 3.8
$row848[0,12] = 1.666
$row848[0,13] = 3.8
$row848[0,14] = 4.75
$row848[0,15] = -0.75
$row848[0,16] = 1.825
$row848[0,17] = 1.975
$row848[0,18] = 2.5
$row848[0,19] = 1.825
$row848[0,20] = 1.975
$row848[0,21] = -1
$row848[0,22] = 2.8
$row848[0,23] = -1
$row848[0,24] = -1
$row848[0,25] = 0.9750000000000001
$row848[0,26] = -1
$row848[0,27] = 0.9750000000000001
$ws.Range("B848:AC848").Value = $row848

$row849 = New-Object "object[,]" 1,28
$row849[0,0] = 6803740
$row849[0,1] = 'Poland I Liga'
$row849[0,2] = 'Poland I Liga'
$row849[0,3] = 45191.54166666666
$row849[0,4] = 'Miedz Legnica'
$row849[0,5] = 'Odra Opole'
$row849[0,6] = 1
$row849[0,7] = 2
$row849[0,8] = 'A'
$row849[0,9] = 1.85
$row849[0,10] = 3.5
$row849[0,11] = 3.75
$row849[0,12] = 1.909
$row849[0,13] = 3.5
$row849[0,14] = 3.5
$row849[0,15] = -0.5
$row849[0,16] = 1.975
$row849[0,17] = 1.825
$row849[0,18] = 2.25
$row849[0,19] = 1.9
$row849[0,20] = 1.9
$row849[0,21] = -1
$row849[0,22] = -1
$row849[0,23] = 2.5
$row849[0,24] = -1
$row849[0,25] = 0.825
$row849[0,26] = 0.8999999999999999
$row849[0,27] = -1
$ws.Range("B849:AC849").Value = $row849

$row909 = New-Object "object[,]" 1,28
$row909[0,0] = 6803778
$row909[0,1] = 'Poland I Liga'
$row909[0,2] = 'Poland I Liga'
$row909[0,3] = 45242.58333333334
$row909[0,4] = 'Podbeskidzie Bielsko Biala'
$row909[0,5] = 'Miedz Legnica'
$row909[0,6] = 0
$row909[0,7] = 0
$row909[0,8] = 'D'
$row909[0,9] = 3.1
$row909[0,10] = 3.4
$row909[0,11] = 2.1
$row909[0,12] = 3
$row909[0,13] = 3.4
$row909[0,14] = 2.15
$row909[0,15] = 0.25
$row909[0,16] = 1.9
$row909[0,17] = 1.95
$row909[0,18] = 2.5
$row909[0,19] = 1.85
$row909[0,20] = 2
$row909[0,21] = -1
$row909[0,22] = 2.4
$row909[0,23] = -1
$row909[0,24] = 0.45
$row909[0,25] = -0.5
$row909[0,26] = -1
$row909[0,27] = 1
$ws.Range("B909:AC909").Value = $row909

$row910 = New-Object "object[,]" 1,28
$row910[0,0] = 6803779
$row910[0,1] = 'Poland I Liga'
$row910[0,2] = 'Poland I Liga'
$row910[0,3] = 45242.58333333334
$row910[0,4] = 'Zaglebie Sosnowiec'
$row910[0,5] = 'Arka Gdynia'
$row910[0,6] = 1
$row910[0,7] = 3
$row910[0,8] = 'A'
$row910[0,9] = 3.3
$row910[0,10] = 3.4
$row910[0,11] = 2
$row910[0,12] = 4.2
$row910[0,13] = 3.5
$row910[0,14] = 1.727
$row910[0,15] = 0.75
$row910[0,16] = 1.825
$row910[0,17] = 2.025
$row910[0,18] = 2.5
$row910[0,19] = 1.9
$row910[0,20] = 1.95
$row910[0,21] = -1
$row910[0,22] = -1
$row910[0,23] = 0.7270000000000001
$row910[0,24] = 1.025
$row910[0,25] = 1
$row910[0,26] = 0.8999999999999999
$row910[0,27] = -1
$ws.Range("B910:AC910").Value = $row910
